$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X3").Value = 686500
$ws.Range("Y3").Value = 151210
$ws.Range("X4").Value = 492200
$ws.Range("Y4").Value = 89405
$ws.Range("X5").Value = 139500
$ws.Range("Y5").Value = 165500
$ws.Range("X6").Value = 157000
$ws.Range("Y6").Value = 174000
$ws.Range("X7").Value = 590500
$ws.Range("Y7").Value = 176325
$ws.Range("X8").Value = 164500
$ws.Range("Y8").Value = 113700
$ws.Range("X9").Value = 289500
$ws.Range("Y9").Value = 347750
$ws.Range("X10").Value = 564600
$ws.Range("Y10").Value = 119000
$ws.Range("X11").Value = 443000
$ws.Range("Y11").Value = 2007660
$ws.Range("F12").Value = 0
$ws.Range("X12").Value = 300500
$ws.Range("X13").Value = 136000
$ws.Range("Y13").Value = 130500
$ws.Range("X14").Value = 211500
$ws.Range("Y14").Value = 78650
$ws.Range("X15").Value = 267500
$ws.Range("Y15").Value = 218750
$ws.Range("X16").Value = 399800
$ws.Range("Y16").Value = 230300
$ws.Range("X17").Value = 381600
$ws.Range("Y17").Value = 91500
$ws.Range("X18").Value = 249500
$ws.Range("Y18").Value = 164500
$ws.Range("X19").Value = 95000
$ws.Range("Y19").Value = 187150
$ws.Range("X20").Value = 156000
$ws.Range("Y20").Value = 335327.39623
$ws.Range("X21").Value = 440500
$ws.Range("Y21").Value = 133200
$ws.Range("X22").Value = 150400
$ws.Range("Y22").Value = 296300
$ws.Range("X23").Value = 752500
$ws.Range("Y23").Value = 100150
$ws.Range("X24").Value = 771000
$ws.Range("Y24").Value = 154750
$ws.Range("X25").Value = 845500
$ws.Range("Y25").Value = 166750
$ws.Range("X26").Value = 638500
$ws.Range("Y26").Value = 1359175
$ws.Range("X27").Value = 273000
$ws.Range("Y27").Value = 313250
$ws.Range("X28").Value = 931000
$ws.Range("Y28").Value = 406525
$ws.Range("X29").Value = 419000
$ws.Range("Y29").Value = 187950
$ws.Range("X30").Value = 375000
$ws.Range("Y30").Value = 360050
$ws.Range("X31").Value = 745250
$ws.Range("Y31").Value = 125900
$ws.Range("X32").Value = 276500
$ws.Range("Y32").Value = 403772.26
$ws.Range("X33").Value = 214150
$ws.Range("Y33").Value = 805850
$ws.Range("X34").Value = 446500
$ws.Range("Y34").Value = 63000
$ws.Range("X35").Value = 186600
$ws.Range("Y35").Value = 323850
$ws.Range("X36").Value = 245500
$ws.Range("Y36").Value = 127600
$ws.Range("X37").Value = 302000
$ws.Range("Y37").Value = 122650
$ws.Range("X38").Value = 481500
$ws.Range("Y38").Value = 176750
$ws.Range("X39").Value = 615250
$ws.Range("Y39").Value = 663600
$ws.Range("X40").Value = 214500
$ws.Range("Y40").Value = 135800
$ws.Range("X41").Value = 209250
$ws.Range("Y41").Value = 111450
$ws.Range("X42").Value = 494250
$ws.Range("Y42").Value = 2530650
$ws.Range("X165").Value = 1028750
$ws.Range("Y165").Value = 2236800
$ws.Range("X166").Value = 470500
$ws.Range("Y166").Value = 347425
$ws.Range("X167").Value = 923250
$ws.Range("Y167").Value = 487000
$ws.Range("X168").Value = 1560000
$ws.Range("Y168").Value = 399450
$ws.Range("X169").Value = 1383000
$ws.Range("Y169").Value = 1147050
$ws.Range("X170").Value = 1809000
$ws.Range("Y170").Value = 343950
$ws.Range("X171").Value = 1222750
$ws.Range("Y171").Value = 283350
$ws.Range("X172").Value = 2078000
$ws.Range("Y172").Value = 2898000
$ws.Range("X173").Value = 676500
$ws.Range("Y173").Value = 205800
$ws.Range("X174").Value = 538250
$ws.Range("Y174").Value = 71600
$ws.Range("X175").Value = 1600500
$ws.Range("Y175").Value = 993225
$ws.Range("X176").Value = 1163750
$ws.Range("Y176").Value = 348600
$ws.Range("X177").Value = 1578000
$ws.Range("Y177").Value = 703650
$ws.Range("X178").Value = 1579750
$ws.Range("Y178").Value = 752800
$ws.Range("X179").Value = 1580500
$ws.Range("Y179").Value = 1848100
$ws.Range("X180").Value = 542250
$ws.Range("Y180").Value = 257150
$ws.Range("X181").Value = 949000
$ws.Range("Y181").Value = 390000
$ws.Range("X182").Value = 1774250
$ws.Range("Y182").Value = 250850
$ws.Range("X183").Value = 985000
$ws.Range("Y183").Value = 409900
$ws.Range("X184").Value = 1147500
$ws.Range("Y184").Value = 516450
$ws.Range("X185").Value = 1131000
$ws.Range("Y185").Value = 371650
$ws.Range("X186").Value = 2120000
$ws.Range("Y186").Value = 1911450
$ws.Range("X187").Value = 444000
$ws.Range("Y187").Value = 939550
$ws.Range("X188").Value = 1008250
$ws.Range("Y188").Value = 126000
$ws.Range("X189").Value = 1591250
$ws.Range("Y189").Value = 330325
$ws.Range("X190").Value = 1170500
$ws.Range("Y190").Value = 1086800
$ws.Range("X191").Value = 1144500
$ws.Range("Y191").Value = 476225
$ws.Range("X192").Value = 1271000
$ws.Range("Y192").Value = 328300
$ws.Range("X193").Value = 1287750
$ws.Range("Y193").Value = 1301450
$ws.Range("X194").Value = 754250
$ws.Range("Y194").Value = 400800
$ws.Range("C195").Value = 77500
$ws.Range("D195").Value = 25500
$ws.Range("F195").Value = 7000
$ws.Range("X195").Value = 1412250
$ws.Range("Y195").Value = 363000
